$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1848
$ws.Range("I28").Value = 478.11765
$ws.Range("K28").Value = 478.11765
$ws.Range("M28").Value = 6.882349999999974

$ws.Range("H92").Value = 2492.6155
$ws.Range("I92").Value = 2390.5
$ws.Range("K92").Value = 2390.5
$ws.Range("M92").Value = -1142.5

$ws.Range("H103").Value = 2006.8572
$ws.Range("J103").Value = 1507.3334
$ws.Range("L103").Value = 4522.0002
$ws.Range("N103").Value = -5694.0002

$ws.Range("H112").Value = 1152.037
$ws.Range("I112").Value = 505.10526
$ws.Range("J112").Value = 2688.5
$ws.Range("K112").Value = 1515.31578
$ws.Range("L112").Value = 8065.5
$ws.Range("M112").Value = -407.3157799999999
$ws.Range("N112").Value = -10281.5

$ws.Range("H138").Value = 1676.8462
$ws.Range("I138").Value = 1527.1818
$ws.Range("J138").Value = 2500
$ws.Range("K138").Value = 4581.5454
$ws.Range("L138").Value = 7500
$ws.Range("M138").Value = 558.4546
$ws.Range("N138").Value = -17780

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 8931538
$ws.Range("I45").Value = 2003.7
$ws.Range("J45").Value = 31255374
$ws.Range("K45").Value = 2003.7
$ws.Range("L45").Value = 31255374
$ws.Range("M45").Value = -1626.7
$ws.Range("N45").Value = -31256128

$ws.Range("H88").Value = 1089.5883
$ws.Range("I88").Value = 956.125
$ws.Range("J88").Value = 1208.2222
$ws.Range("K88").Value = 956.125
$ws.Range("L88").Value = 1208.2222
$ws.Range("M88").Value = -550.125
$ws.Range("N88").Value = -2020.2222

$ws.Range("H91").Value = 1089.5883
$ws.Range("I91").Value = 956.125
$ws.Range("J91").Value = 1208.2222
$ws.Range("K91").Value = 956.125
$ws.Range("L91").Value = 1208.2222
$ws.Range("M91").Value = 447.875
$ws.Range("N91").Value = -4016.2222

$ws.Range("H132").Value = 2315.1724
$ws.Range("I132").Value = 2241.2
$ws.Range("J132").Value = 2479.5557
$ws.Range("K132").Value = 6723.599999999999
$ws.Range("L132").Value = 7438.6671
$ws.Range("M132").Value = -4193.599999999999
$ws.Range("N132").Value = -12498.6671

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 115555.555
$ws.Range("I20").Value = 115555.555
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 115555.555
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -115308.555
$ws.Range("N20").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 349.16666
$ws.Range("J22").Value = 1000
$ws.Range("L22").Value = 1000
$ws.Range("N22").Value = -1700

$ws.Range("H86").Value = 3577418.2
$ws.Range("I86").Value = 5956613.5
$ws.Range("K86").Value = 5956613.5
$ws.Range("M86").Value = -5955490.5

$ws.Range("H89").Value = 3577418.2
$ws.Range("I89").Value = 5956613.5
$ws.Range("K89").Value = 29783067.5
$ws.Range("M89").Value = -29777451.5

$ws.Range("H94").Value = 2000
$ws.Range("I94").Value = 2000
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 2000
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -1549
$ws.Range("N94").ClearContents()

$ws.Range("H132").Value = 1721133.4
$ws.Range("I132").Value = 2274598.5
$ws.Range("K132").Value = 6823795.5
$ws.Range("M132").Value = -6821265.5

$ws.Range("H133").Value = 24326
$ws.Range("J133").Value = 24326
$ws.Range("L133").Value = 24326
$ws.Range("N133").Value = -29386

$ws.Range("H134").Value = 3451858.2
$ws.Range("I134").Value = 4467251.5
$ws.Range("K134").Value = 13401754.5
$ws.Range("M134").Value = -13399219.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 732.65717
$ws.Range("J5").Value = 943.2308
$ws.Range("L5").Value = 2829.6924
$ws.Range("N5").Value = -3053.6924

$ws.Range("H33").Value = 263.5
$ws.Range("I33").Value = 222.72728
$ws.Range("J33").Value = 413
$ws.Range("K33").Value = 1336.36368
$ws.Range("L33").Value = 2478
$ws.Range("M33").Value = -1053.36368
$ws.Range("N33").Value = -3044

$ws.Range("H68").Value = 102681.9
$ws.Range("J68").Value = 114054
$ws.Range("L68").Value = 342162
$ws.Range("N68").Value = -343784

$ws.Range("H71").Value = 102681.9
$ws.Range("J71").Value = 114054
$ws.Range("L71").Value = 1026486
$ws.Range("N71").Value = -1034598

$ws.Range("H119").Value = 2837.5
$ws.Range("J119").Value = 4000
$ws.Range("L119").Value = 12000
$ws.Range("N119").Value = -21676

$ws.Range("H135").Value = 732.65717
$ws.Range("J135").Value = 943.2308
$ws.Range("L135").Value = 8489.0772
$ws.Range("N135").Value = -13559.0772

$ws.Range("H141").Value = 1750
$ws.Range("I141").Value = 1750
$ws.Range("K141").Value = 5250
$ws.Range("M141").Value = -70

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3972.2
$ws.Range("I132").Value = 4121.923
$ws.Range("J132").Value = 2999
$ws.Range("K132").Value = 12365.769
$ws.Range("L132").Value = 8997
$ws.Range("M132").Value = -9835.769
$ws.Range("N132").Value = -14057

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1233.2632
$ws.Range("I22").Value = 1381.8334
$ws.Range("J22").Value = 978.5714
$ws.Range("K22").Value = 1381.8334
$ws.Range("L22").Value = 978.5714
$ws.Range("M22").Value = -1086.8334
$ws.Range("N22").Value = -1568.5714

$ws.Range("H27").Value = 1233.2632
$ws.Range("I27").Value = 1381.8334
$ws.Range("J27").Value = 978.5714
$ws.Range("K27").Value = 1381.8334
$ws.Range("L27").Value = 978.5714
$ws.Range("M27").Value = -1274.8334
$ws.Range("N27").Value = -1192.5714

$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()

$ws.Range("H100").Value = 4320.9565
$ws.Range("I100").Value = 3780.1428
$ws.Range("K100").Value = 3780.1428
$ws.Range("M100").Value = -3239.1428

$ws.Range("H136").Value = 2998.5
$ws.Range("I136").Value = 2500
$ws.Range("J136").Value = 3330.8333
$ws.Range("K136").Value = 7500
$ws.Range("L136").Value = 9992.499899999999
$ws.Range("M136").Value = -4950
$ws.Range("N136").Value = -15092.4999

$ws.Range("H141").Value = 72851.60000000001
$ws.Range("J141").Value = 72851.60000000001
$ws.Range("L141").Value = 72851.60000000001
$ws.Range("N141").Value = -83211.60000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 30023.334
$ws.Range("I51").Value = 30023.334
$ws.Range("K51").Value = 30023.334
$ws.Range("M51").Value = -29513.334

$ws.Range("H100").Value = 3106321
$ws.Range("I100").Value = 3969005.5
$ws.Range("K100").Value = 7938011
$ws.Range("M100").Value = -7937470

$ws.Range("H107").Value = 2143.5
$ws.Range("I107").Value = 1466.0834
$ws.Range("J107").Value = 4175.75
$ws.Range("K107").Value = 4398.2502
$ws.Range("L107").Value = 12527.25
$ws.Range("M107").Value = -2478.2502
$ws.Range("N107").Value = -16367.25

$ws.Range("H132").Value = 2092.318
$ws.Range("I132").Value = 1936.8462
$ws.Range("K132").Value = 5810.5386
$ws.Range("M132").Value = -3280.5386
